# New predicted prices and compared returns.
# Updates the "Return_with_prediction" (G) and "return_pct_change" (H)
# columns for rows 2-85, plus the single "mean_return_pct_change" (I2) cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @(0.06854349825647522,0.06685146085328561,-0.01294739277674126,-0.00805072599712036,-0.01860918511944517,-0.007200376581823696,-0.006208695976173156,0.004615173564844266,-0.05921775018514164,-0.06954880205389437,-0.3777130478920371,-0.3974577424205049,-0.02013216919123732,-0.0140782835910078,0.1421372740448495,0.1419450036003737,0.1214636032858463,0.1280493617729934,0.08583424832767612,0.08750101757515669,-0.09776778619565953,-0.09996788445127837,0.1622172585708786,0.1686872436201788,0.08861000624943989,0.08724696713513973,-0.128692396270753,-0.1427165399519818,0.04987456996619302,0.04719338620083045,0.1126956322543052,0.1189988281279437,-0.008899927039406726,-0.01336337996427005,0.0332671601989017,0.03917196530920936,0.1056205784144961,0.09849587206358534,0.03085941075929714,0.03739420001845033,0.1142063130845872,0.1283173601802984,0.04154277884281783,0.03205172709402848,0.0625297954706884,0.06170374293006865,0.04445528198183851,0.05131862660928519,0.02929800480279795,0.03024472037892732,-0.08203549445597172,-0.08241151936368532,0.05486635106409896,0.05442158487578165,0.05097621793781661,0.04215395513469582,0.0535991321761355,0.061031896618933,0.02189470446886582,0.03473810853425006,0.06489112235944754,0.06113028721643708,0.02488332273580799,0.03489841065810598,0.07381117535184126,0.08420562091045028,-0.0216296105633652,-0.01279102276884447,0.07633135185111055,0.07792875535015219,-0.1355361960702551,-0.1469157655986069,0.1449383962768462,0.1521327215339489,-0.005421420542721624,-0.003564192275308726,0.08963188686374464,0.0945267997282101,-0.2307771327044859,-0.2102747674148371,0.1766910187055766,0.1852974508877285,0.1091612763950627,0.107052983348096)
$hValues = @(6.595883129145573,18.93787868357008,36.25619031462482,29.57315274768672,-66.24906637884965,-27.1312181143994,-9.131975639318677,184.0169433930521,6.03601802054687,-8.537641576057208,4.280336388448042,-1.361459035477968,-145.6861333888471,68.93491822385735,3.939772229358073,1.774111820730376,3.225506305601552,-0.4814603200043607,-3.268604691598187,0.4707629821675995,-4.583302872704445,1.456092867198455,0.6991058805038493,-1.118719865288308,-2.254499576005931,1.495107895281908,6.501025662900731,-2.086394403784631,-4.115530755665511,7.699747043314642,3.672084834344359,-4.106818491361187,43.00907706590569,20.15803492658402,-9.518747614284074,9.769520790917436,5.30232385932141,1.12510978077343,-8.395395724047278,16.06356234946749,-5.540791157133192,0.4162547314318349,4.729767037628193,2.854001958011117,10.44220113771262,5.170524315650755,-9.748887018542678,12.62050391271533,10.61492750571839,7.955676698407597,5.623495049993994,-2.736750658619901,9.674131914210516,-3.310775556441291,3.124142723693939,10.97825323117905,-6.968647841529603,7.048507423061183,-20.26691604945054,30.11587319196584,3.913529932044643,-4.308495818313278,-10.2969356983362,-1.49275268593271,-4.984355372559917,6.7639551227707,0.5338168413236078,33.17583121987436,6.036586700223745,-1.878951140427929,11.78411325217087,4.024018886304042,-3.643682333484324,1.12025720639251,-423.0700520626725,-61.42995579324317,-0.3627944644389623,-2.450496960200792,-6.627335230460615,1.332206920941838,5.41312297093401,5.265338253207235,2.856485474985425,2.385429326441588)

$startRow = 2
for ($i = 0; $i -lt $gValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $gValues[$i]
    $ws.Cells.Item($row, 8).Value = $hValues[$i]
}

$ws.Cells.Item(2, 9).Value = -2.145003858990287
